$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Bec: update Y1, Y2, X1, X2
$ws.Range("B5").Value = 896
$ws.Range("C5").Value = 1182
$ws.Range("D5").Value = 1161
$ws.Range("E5").Value = 1453

# Row 6 - SideOdtCamera: update Y1, Y2, X1, X2 and add SubRoiCenterSize
$ws.Range("B6").Value = 436
$ws.Range("C6").Value = 562
$ws.Range("D6").Value = 825
$ws.Range("E6").Value = 963
$ws.Range("I6").Value = "[]"

# Row 7 - NiLattice: update Y1, Y2, X1, X2
$ws.Range("B7").Value = 781
$ws.Range("C7").Value = 1393
$ws.Range("D7").Value = 1233
$ws.Range("E7").Value = 1369

# Row 20 - BecCameraSBB: update Y1, Y2, X1, X2
$ws.Range("B20").Value = 447
$ws.Range("C20").Value = 735
$ws.Range("D20").Value = 921
$ws.Range("E20").Value = 1203
